# Update countries & provincias Spain
# Applies the data refresh + shared-string reordering captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 13 de Julio de 2020 a las 22:22'
$ws.Cells.Item(4, 2).Value = 3464004
$ws.Cells.Item(4, 3).Value = 50009
$ws.Cells.Item(4, 4).Value = 1537766
$ws.Cells.Item(4, 5).Value = 1788150
$ws.Cells.Item(4, 7).Value = 306
$ws.Cells.Item(4, 8).Value = 138088
$ws.Cells.Item(6, 2).Value = 907645
$ws.Cells.Item(6, 3).Value = 28179
$ws.Cells.Item(6, 4).Value = 572112
$ws.Cells.Item(6, 5).Value = 311806
$ws.Cells.Item(23, 4).Value = 71648
$ws.Cells.Item(23, 5).Value = 27372
$ws.Cells.Item(46, 2).Value = 40632
$ws.Cells.Item(46, 3).Value = 1962
$ws.Cells.Item(46, 4).Value = 19395
$ws.Cells.Item(46, 5).Value = 20872
$ws.Cells.Item(67, 1).Value = 'Uzbekistan'
$ws.Cells.Item(67, 2).Value = 13591
$ws.Cells.Item(67, 3).Value = 594
$ws.Cells.Item(67, 4).Value = 8030
$ws.Cells.Item(67, 5).Value = 5497
$ws.Cells.Item(67, 7).Value = 4
$ws.Cells.Item(67, 8).Value = 64
$ws.Cells.Item(68, 1).Value = 'Corea del Sur'
$ws.Cells.Item(68, 2).Value = 13479
$ws.Cells.Item(68, 3).Value = 62
$ws.Cells.Item(68, 4).Value = 12204
$ws.Cells.Item(68, 5).Value = 986
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 289
$ws.Cells.Item(73, 1).Value = 'Sudan'
$ws.Cells.Item(73, 2).Value = 10316
$ws.Cells.Item(73, 3).Value = 66
$ws.Cells.Item(73, 4).Value = 5403
$ws.Cells.Item(73, 5).Value = 4256
$ws.Cells.Item(73, 7).Value = 7
$ws.Cells.Item(73, 8).Value = 657
$ws.Cells.Item(74, 1).Value = 'Kenia'
$ws.Cells.Item(74, 2).Value = 10294
$ws.Cells.Item(74, 3).Value = 189
$ws.Cells.Item(74, 4).Value = 2946
$ws.Cells.Item(74, 5).Value = 7151
$ws.Cells.Item(74, 7).Value = 12
$ws.Cells.Item(74, 8).Value = 197
$ws.Cells.Item(83, 1).Value = 'Costa Rica'
$ws.Cells.Item(83, 2).Value = 8036
$ws.Cells.Item(83, 3).Value = 440
$ws.Cells.Item(83, 4).Value = 2304
$ws.Cells.Item(83, 5).Value = 5701
$ws.Cells.Item(83, 8).Value = 31
$ws.Cells.Item(84, 1).Value = 'Etiopia'
$ws.Cells.Item(84, 2).Value = 7766
$ws.Cells.Item(84, 3).Value = 206
$ws.Cells.Item(84, 4).Value = 2430
$ws.Cells.Item(84, 5).Value = 5208
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = 128
$ws.Cells.Item(94, 2).Value = 5446
$ws.Cells.Item(94, 3).Value = 91
$ws.Cells.Item(94, 4).Value = 2501
$ws.Cells.Item(94, 5).Value = 2798
$ws.Cells.Item(98, 2).Value = 4321
$ws.Cells.Item(98, 3).Value = 33
$ws.Cells.Item(98, 4).Value = 1208
$ws.Cells.Item(98, 5).Value = 3060
$ws.Cells.Item(106, 2).Value = 2980
$ws.Cells.Item(106, 3).Value = 32
$ws.Cells.Item(106, 4).Value = 1293
$ws.Cells.Item(106, 5).Value = 1665
$ws.Cells.Item(109, 2).Value = 2724
$ws.Cells.Item(109, 3).Value = 13
$ws.Cells.Item(109, 5).Value = 207
$ws.Cells.Item(111, 1).Value = 'Malaui'
$ws.Cells.Item(111, 2).Value = 2430
$ws.Cells.Item(111, 3).Value = 66
$ws.Cells.Item(111, 4).Value = 747
$ws.Cells.Item(111, 5).Value = 1644
$ws.Cells.Item(111, 7).Value = 1
$ws.Cells.Item(111, 8).Value = 39
$ws.Cells.Item(112, 1).Value = 'Cuba'
$ws.Cells.Item(112, 2).Value = 2428
$ws.Cells.Item(112, 3).Value = 2
$ws.Cells.Item(112, 4).Value = 2268
$ws.Cells.Item(112, 5).Value = 73
$ws.Cells.Item(112, 8).Value = 87
$ws.Cells.Item(113, 1).Value = 'Libano'
$ws.Cells.Item(113, 2).Value = 2419
$ws.Cells.Item(113, 3).Value = 85
$ws.Cells.Item(113, 4).Value = 1423
$ws.Cells.Item(113, 5).Value = 960
$ws.Cells.Item(113, 8).Value = 36
$ws.Cells.Item(114, 1).Value = 'Mali'
$ws.Cells.Item(114, 2).Value = 2412
$ws.Cells.Item(114, 3).Value = 1
$ws.Cells.Item(114, 4).Value = 1730
$ws.Cells.Item(114, 5).Value = 561
$ws.Cells.Item(114, 8).Value = 121
$ws.Cells.Item(146, 4).Value = 792
$ws.Cells.Item(146, 5).Value = 13
$ws.Cells.Item(209, 1).Value = 'Groenlandia'
$ws.Cells.Item(210, 1).Value = 'Islas Malvinas'
